# fix: Member table sync, My Account display, QR dep
#
# 1. Member table sync: backfill the mobile number for member MEM-003
#    (Ganeshan) on his two existing, still-open loans (rows 3 and 9) so
#    the sheet stays in sync with his profile/member table.
# 2. QR dep: a new loan was checked out and returned via the QR
#    check-in/out flow - GDL-008 borrowed and returned by Ganeshan
#    (MEM-003) - recorded as a new row with a freshly generated
#    transaction id.
# 3. My Account display: this also fixes the "My Account" transaction
#    history, which was missing this QR-driven loan/return.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add missing mobile numbers for Ganeshan (MEM-003) on rows 3 and 9
$ws.Range("G3").Value = 1234567899
$ws.Range("G9").Value = 1234567899

# 2) New row 31: Ganeshan's GDL-008 loan, checked out/in via QR, with a
#    freshly generated transaction id
$ws.Range("A31").Value = "TX-20260117201635"
$ws.Range("B31").Value = "GDL-008"
$ws.Range("C31").Value = "அக்னி சிறகுகள் (2 COPIES)"
$ws.Range("D31").Value = "MEM-003"
$ws.Range("E31").Value = "Ganeshan"
$ws.Range("G31").Value = 1234567899
# leading apostrophe forces these date-shaped values to stay plain text,
# matching how borrow_date/return_date are stored elsewhere in the sheet
$ws.Range("H31").Value = "'2026-01-17"
$ws.Range("I31").Value = "'2026-01-17"
$ws.Range("J31").Value = "RETURNED"
